# Update LR-pair NATMI stats (Tgfb3-Acvrl1) with recomputed TPM-based values.
# Only numeric value cells in columns G,H,I,J (ligand stats), M,N,O,P (receptor
# stats) and Q,R,S,T (edge stats) for rows 2-10 change; everything else is
# left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.103903333333333
$ws.Range("H2").Value = 3.31171
$ws.Range("I2").Value = 0.02393122995918198
$ws.Range("J2").Value = 0.02393122995918198
$ws.Range("M2").Value = 19.92674333333333
$ws.Range("N2").Value = 59.78023
$ws.Range("O2").Value = 0.3447897148135736
$ws.Range("P2").Value = 0.3447897148135735
$ws.Range("Q2").Value = 21.99719838814445
$ws.Range("R2").Value = 197.9747854933
$ws.Range("S2").Value = 0.008251241952764402
$ws.Range("T2").Value = 0.008251241952764402
$ws.Range("G3").Value = 1.103903333333333
$ws.Range("H3").Value = 3.31171
$ws.Range("I3").Value = 0.02393122995918198
$ws.Range("J3").Value = 0.02393122995918198
$ws.Range("O3").Value = 0.5793221821339875
$ws.Range("P3").Value = 0.5793221821339873
$ws.Range("Q3").Value = 36.96010763530001
$ws.Range("R3").Value = 332.6409687177
$ws.Range("S3").Value = 0.01386389236110356
$ws.Range("T3").Value = 0.01386389236110356
$ws.Range("G4").Value = 1.103903333333333
$ws.Range("H4").Value = 3.31171
$ws.Range("I4").Value = 0.02393122995918198
$ws.Range("J4").Value = 0.02393122995918198
$ws.Range("M4").Value = 4.385869666666667
$ws.Range("N4").Value = 13.157609
$ws.Range("O4").Value = 0.07588810305243907
$ws.Range("P4").Value = 0.07588810305243905
$ws.Range("Q4").Value = 4.84157614459889
$ws.Range("R4").Value = 43.57418530139
$ws.Range("S4").Value = 0.001816095645314019
$ws.Range("T4").Value = 0.001816095645314019
$ws.Range("H5").Value = 70.73212899999999
$ws.Range("I5").Value = 0.5111277390235027
$ws.Range("J5").Value = 0.5111277390235027
$ws.Range("M5").Value = 19.92674333333333
$ws.Range("N5").Value = 59.78023
$ws.Range("O5").Value = 0.3447897148135736
$ws.Range("P5").Value = 0.3447897148135735
$ws.Range("Q5").Value = 469.8203266677411
$ws.Range("R5").Value = 4228.382940009669
$ws.Range("S5").Value = 0.1762315873712202
$ws.Range("T5").Value = 0.1762315873712201
$ws.Range("H6").Value = 70.73212899999999
$ws.Range("I6").Value = 0.5111277390235027
$ws.Range("J6").Value = 0.5111277390235027
$ws.Range("O6").Value = 0.5793221821339875
$ws.Range("P6").Value = 0.5793221821339873
$ws.Range("Q6").Value = 789.40097445547
$ws.Range("R6").Value = 7104.608770099229
$ws.Range("S6").Value = 0.2961076371203069
$ws.Range("T6").Value = 0.2961076371203068
$ws.Range("H7").Value = 70.73212899999999
$ws.Range("I7").Value = 0.5111277390235027
$ws.Range("J7").Value = 0.5111277390235027
$ws.Range("M7").Value = 4.385869666666667
$ws.Range("N7").Value = 13.157609
$ws.Range("O7").Value = 0.07588810305243907
$ws.Range("P7").Value = 0.07588810305243905
$ws.Range("Q7").Value = 103.4072996799512
$ws.Range("R7").Value = 930.6656971195608
$ws.Range("S7").Value = 0.03878851453197576
$ws.Range("T7").Value = 0.03878851453197575
$ws.Range("G8").Value = 21.446869
$ws.Range("H8").Value = 64.34060699999999
$ws.Range("I8").Value = 0.4649410310173153
$ws.Range("J8").Value = 0.4649410310173154
$ws.Range("M8").Value = 19.92674333333333
$ws.Range("N8").Value = 59.78023
$ws.Range("O8").Value = 0.3447897148135736
$ws.Range("P8").Value = 0.3447897148135735
$ws.Range("Q8").Value = 427.3662538666233
$ws.Range("R8").Value = 3846.29628479961
$ws.Range("S8").Value = 0.160306885489589
$ws.Range("T8").Value = 0.160306885489589
$ws.Range("G9").Value = 21.446869
$ws.Range("H9").Value = 64.34060699999999
$ws.Range("I9").Value = 0.4649410310173153
$ws.Range("J9").Value = 0.4649410310173154
$ws.Range("O9").Value = 0.5793221821339875
$ws.Range("P9").Value = 0.5793221821339873
$ws.Range("Q9").Value = 718.0688405810099
$ws.Range("R9").Value = 6462.619565229089
$ws.Range("S9").Value = 0.2693506526525771
$ws.Range("T9").Value = 0.269350652652577
$ws.Range("G10").Value = 21.446869
$ws.Range("H10").Value = 64.34060699999999
$ws.Range("I10").Value = 0.4649410310173153
$ws.Range("J10").Value = 0.4649410310173154
$ws.Range("M10").Value = 4.385869666666667
$ws.Range("N10").Value = 13.157609
$ws.Range("O10").Value = 0.07588810305243907
$ws.Range("P10").Value = 0.07588810305243905
$ws.Range("Q10").Value = 94.06317219207367
$ws.Range("R10").Value = 846.5685497286629
$ws.Range("S10").Value = 0.03528349287514929
$ws.Range("T10").Value = 0.03528349287514929
